$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H
$ws.Range("H1").Value = "FRA"

# FRA classification values for rows 2-23 (Y/N/U)
$values = @("Y","Y","Y","N","N","Y","Y","Y","U","N","N","U","Y","Y","Y","Y","Y","Y","Y","N","N","N")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

$ws.Cells.Item(24, 8).Value = ""

$ws.Range("A11").Select()
$ws.Range("H24").Select()
